$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 813, shifting existing rows 813-854 down to 814-855.
$ws.Rows(813).Insert()

# Populate the newly inserted row with the new data point for 2026/02/15.
# Force column A to remain literal text (matches the rest of the date
# column) instead of letting Excel auto-convert the "yyyy/mm/dd"-looking
# string into a real date serial value, then restore the default "Normal"
# style so no stray number-format style lingers on the cell.
$ws.Cells.Item(813, 1).NumberFormat = "@"
$ws.Cells.Item(813, 1).Value = "2026/02/15"
$ws.Cells.Item(813, 1).Style = "Normal"
$ws.Cells.Item(813, 2).Value = "日"
$ws.Cells.Item(813, 3).Value = 8
$ws.Cells.Item(813, 4).Value = 201
